$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values (row 1), columns O and P -> 14, 15
$ws.Range("O1").Value = 14
$ws.Range("P1").Value = 15

# Copy style from N1 (header style) to O1 and P1
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("O2").Value = -0.8366367068726722
$ws.Range("P2").Value = -0.6126118682264818

$ws.Range("O3").Value = -0.4657240759688368
$ws.Range("P3").Value = -0.4001321118272729

$ws.Range("O4").Value = 0.03754329094432164
$ws.Range("P4").Value = -0.008940944522873105

$ws.Range("O5").Value = 0.4166668696793512
$ws.Range("P5").Value = 0.3768613040501999

$ws.Range("O6").Value = -0.3006722170473243
$ws.Range("P6").Value = -0.2919641625736087

$ws.Range("O7").Value = -0.1577104617054219
$ws.Range("P7").Value = -0.1574877845208207

$ws.Range("O8").Value = -0.442120638476823
$ws.Range("P8").Value = -0.4342224082364093

$ws.Range("O9").Value = 0.003556211750688697
$ws.Range("P9").Value = 0.002969031900047409

$ws.Range("O10").Value = 0.006654235906479544
$ws.Range("P10").Value = 0.007351302360603446

$ws.Range("O11").Value = 0.009029670688592699
$ws.Range("P11").Value = 0.007944550950589291
